# Edit applies two changes to the "Problem, issues and concerns" numbered
# list (numId=4):
#   1. Insert a new, empty list item between the "There are two circuit
#      design schemes..." item and the "to design the output impandance..."
#      item.
#   2. Fix the typo in "May need a bandpass near output" (becomes
#      "May a bandp need ass near output") and remove the following,
#      now-redundant "Think how" list item entirely.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $exactText) {
    $match = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq ($exactText + [char]13)) {
            $match = $p
        }
    }
    return $match
}

# --- 1. Insert a new empty ListParagraph list item -------------------------
$schemesPara = Find-ParagraphByText $d "There are two circuit design schemes we confirmed and it is hard to choose which one to use."
$endPos = $schemesPara.Range.End
# Re-create the insertion point as a brand-new document Range (rather than a
# collapsed duplicate of the paragraph's own Range) - InsertXML only behaves
# correctly (inserting a standalone paragraph break) when given a "fresh"
# Range object at that offset.
$insertionPoint = $d.Range($endPos, $endPos)

# Build a WordOpenXML "pkg:package" fragment containing two <w:p> elements:
# the boundary between them becomes the new paragraph break inserted at
# $insertionPoint. The first <w:p> supplies the pPr/rPr for the brand new
# (empty) list paragraph; the second <w:p> carries a uniquely-named
# placeholder run so the following paragraph's own runs are not silently
# merged/eaten by the insertion.
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t>ZZ_PLACEHOLDER_ZZ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xmlFrag) | Out-Null

# Remove the placeholder run we had to plant to keep the following
# paragraph's own content intact.
$d.Content.Find.Execute("ZZ_PLACEHOLDER_ZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- 2. Fix the typo, then drop the obsolete "Think how" item --------------
$d.Content.Find.Execute("May need a bandpass near output", $true, $false, $false, $false, $false, $true, 1, $false, "May a bandp need ass near output", 2) | Out-Null

$thinkHowPara = Find-ParagraphByText $d "Think how"
if ($thinkHowPara -ne $null) {
    $thinkHowPara.Range.Delete() | Out-Null
}
